$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, border, centered) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Re-assert the header text values (PasteSpecial only carries formats)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Values for I2:I70 and J2:J70 taken from the diff (row 2 .. row 70)
$iVals = @(9,9,9,9,7,9,9,7,8,9,9,7,9,8,8,9,8,7,8,9,7,9,9,8,9,6,7,9,7,7,7,8,8,8,8,7,8,10,6,7,9,9,9,8,4,9,9,6,6,7,5,7,5,5,1,5,9,8,8,5,6,7,7,6,6,9,7,6,4)
$jVals = @(9,9,9,9,8,9,9,7,8,9,9,7,9,8,8,9,8,9,9,9,7,9,9,9,9,6,7,9,7,7,7,8,8,8,8,7,8,10,6,7,9,9,9,8,4,9,9,6,7,7,6,7,6,6,2,6,9,8,8,6,7,7,7,6,6,9,7,6,4)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
